# average over return periods works well with single and multi hazard data
#
# Adds a new "Discount rate" / "rho" column (Y) to the compiled data sheet,
# filled with a constant discount rate of 0.05 for every province data row,
# and normalizes the "Province" header label to lowercase "province".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column Y: Discount rate / rho -------------------------------------

# Header row (row 1) - long description, matches style of the other header
# cells (bold, centered, bordered).
$ws.Range("X1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Y1").Value = "Discount rate"

# Variable-name row (row 2) - short code, same header styling.
$ws.Range("X2").Copy()
$ws.Range("Y2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Y2").Value = "rho"

# Data rows (4 through 83) - constant discount rate value. Row 3 is the
# section/label row (only column A is populated there), so it is skipped.
for ($r = 4; $r -le 83; $r++) {
    $ws.Cells.Item($r, 25).Value = 0.05
}

# --- Rename "Province" label to lowercase "province" -----------------------

$ws.Range("A3").Value = "province"
